# Add another line item (row 3) to the Expenses sheet, mirroring the
# structure/style of the existing row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row at position 3. Excel's default "insert" behavior
# copies the cell formatting (style) of the row above (row 2) into the new
# row, column by column, which is exactly what we need since row 3 reuses
# the same per-column styles as row 2 in the target file.
$ws.Rows(3).Insert()

# Populate the new row's values.
$ws.Range("A3").Value = "Power"
$ws.Range("B3").Value = "Computers"
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 23000
$ws.Range("E3").Value = "DKK"
$ws.Range("F3").Value = 44714.51493055555
$ws.Range("G3").Value = "https://www.power.dk/computere-og-tablets/computere/baerbar-pc/hp-14s-fq2473no-14-baerbar-pc/p-1846428/"
$ws.Range("H3").Value = "HP Bærbar"
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = 5.52
$ws.Range("K3").Value = $true
$ws.Range("L3").Value = 44987.002916666665
$ws.Range("M3").Value = 44584
$ws.Range("N3").Value = 25569.48130787037
$ws.Range("O3").Value = 5.52

# Wire up the hyperlinks for the Receipt (G) and URL (H) columns, both
# pointing at the product page, matching the pattern used in row 2.
$ws.Hyperlinks.Add($ws.Range("G3"), "https://www.power.dk/computere-og-tablets/computere/baerbar-pc/hp-14s-fq2473no-14-baerbar-pc/p-1846428/")
$ws.Hyperlinks.Add($ws.Range("H3"), "https://www.power.dk/computere-og-tablets/computere/baerbar-pc/hp-14s-fq2473no-14-baerbar-pc/p-1846428/")
